$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(99,8).Value = 1576
$ws.Cells.Item(99,9).Value = 1576
$ws.Cells.Item(99,11).Value = 4728
$ws.Cells.Item(99,13).Value = -3230

$ws.Cells.Item(100,8).Value = 1626.24
$ws.Cells.Item(100,9).Value = 885.0769
$ws.Cells.Item(100,10).Value = 2429.1667
$ws.Cells.Item(100,11).Value = 885.0769
$ws.Cells.Item(100,12).Value = 2429.1667
$ws.Cells.Item(100,13).Value = -344.0769
$ws.Cells.Item(100,14).Value = -3511.1667

$ws.Cells.Item(101,8).Value = 224.75
$ws.Cells.Item(101,9).Value = 224.75
$ws.Cells.Item(101,10).Value = 0
$ws.Cells.Item(101,11).Value = 674.25
$ws.Cells.Item(101,12).Value = 0
$ws.Cells.Item(101,13).Value = 947.75
$ws.Cells.Item(101,14).ClearContents()

$ws.Cells.Item(141,8).Value = 5898.4614
$ws.Cells.Item(141,9).Value = 3880
$ws.Cells.Item(141,10).Value = 17000
$ws.Cells.Item(141,11).Value = 11640
$ws.Cells.Item(141,12).Value = 51000
$ws.Cells.Item(141,13).Value = -6460
$ws.Cells.Item(141,14).Value = -61360

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(101,8).Value = 20000
$ws.Cells.Item(101,10).Value = 20000
$ws.Cells.Item(101,12).Value = 20000
$ws.Cells.Item(101,14).Value = -26490

$ws.Cells.Item(102,8).Value = 1558.7059
$ws.Cells.Item(102,9).Value = 1538.3077
$ws.Cells.Item(102,10).Value = 1625
$ws.Cells.Item(102,11).Value = 1538.3077
$ws.Cells.Item(102,12).Value = 1625
$ws.Cells.Item(102,13).Value = 83.69229999999993
$ws.Cells.Item(102,14).Value = -4869

$ws.Cells.Item(110,8).Value = 3054.2
$ws.Cells.Item(110,9).Value = 3515.6667
$ws.Cells.Item(110,11).Value = 3515.6667
$ws.Cells.Item(110,13).Value = -1470.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105,8).Value = 27490
$ws.Cells.Item(105,10).Value = 27490
$ws.Cells.Item(105,12).Value = 27490
$ws.Cells.Item(105,14).Value = -30984

$ws.Cells.Item(107,8).Value = 1218.2963
$ws.Cells.Item(107,9).Value = 986.8
$ws.Cells.Item(107,11).Value = 986.8
$ws.Cells.Item(107,13).Value = 933.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58,8).Value = 3266.9167
$ws.Cells.Item(58,9).Value = 910.3077
$ws.Cells.Item(58,10).Value = 6052
$ws.Cells.Item(58,11).Value = 910.3077
$ws.Cells.Item(58,12).Value = 6052
$ws.Cells.Item(58,13).Value = -707.3077
$ws.Cells.Item(58,14).Value = -6458

$ws.Cells.Item(62,8).Value = 41669380
$ws.Cells.Item(62,9).Value = 2413
$ws.Cells.Item(62,10).Value = 71431500
$ws.Cells.Item(62,11).Value = 2413
$ws.Cells.Item(62,12).Value = 71431500
$ws.Cells.Item(62,13).Value = -1789
$ws.Cells.Item(62,14).Value = -71432748

$ws.Cells.Item(65,8).Value = 41669380
$ws.Cells.Item(65,9).Value = 2413
$ws.Cells.Item(65,10).Value = 71431500
$ws.Cells.Item(65,11).Value = 12065
$ws.Cells.Item(65,12).Value = 357157500
$ws.Cells.Item(65,13).Value = -8945
$ws.Cells.Item(65,14).Value = -357163740

$ws.Cells.Item(132,8).Value = 2217.04
$ws.Cells.Item(132,9).Value = 1570.5625
$ws.Cells.Item(132,10).Value = 3366.3333
$ws.Cells.Item(132,11).Value = 4711.6875
$ws.Cells.Item(132,12).Value = 10098.9999
$ws.Cells.Item(132,13).Value = -2181.6875
$ws.Cells.Item(132,14).Value = -15158.9999

$ws.Cells.Item(136,8).Value = 3266.9167
$ws.Cells.Item(136,9).Value = 910.3077
$ws.Cells.Item(136,10).Value = 6052
$ws.Cells.Item(136,11).Value = 2730.9231
$ws.Cells.Item(136,12).Value = 18156
$ws.Cells.Item(136,13).Value = -180.9231
$ws.Cells.Item(136,14).Value = -23256

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(36,8).Value = 394.875
$ws.Cells.Item(36,9).Value = 394.875
$ws.Cells.Item(36,10).Value = 0
$ws.Cells.Item(36,11).Value = 1184.625
$ws.Cells.Item(36,12).Value = 0
$ws.Cells.Item(36,13).Value = -1015.625
$ws.Cells.Item(36,14).ClearContents()

$ws.Cells.Item(58,8).Value = 9000
$ws.Cells.Item(58,9).Value = 9000
$ws.Cells.Item(58,10).Value = 9000
$ws.Cells.Item(58,11).Value = 27000
$ws.Cells.Item(58,12).Value = 27000
$ws.Cells.Item(58,13).Value = -26872
$ws.Cells.Item(58,14).Value = -27256

$ws.Cells.Item(131,8).Value = 996029.6
$ws.Cells.Item(131,10).Value = 1108.0358
$ws.Cells.Item(131,12).Value = 3324.1074
$ws.Cells.Item(131,14).Value = -13404.1074

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(20,8).Value = 4006
$ws.Cells.Item(20,10).Value = 4006
$ws.Cells.Item(20,12).Value = 4006
$ws.Cells.Item(20,14).Value = -4496

$ws.Cells.Item(43,8).Value = 919.8182
$ws.Cells.Item(43,9).Value = 919.8182
$ws.Cells.Item(43,11).Value = 919.8182
$ws.Cells.Item(43,13).Value = -768.8182

$ws.Cells.Item(46,8).Value = 6000
$ws.Cells.Item(46,9).Value = 6000
$ws.Cells.Item(46,10).Value = 0
$ws.Cells.Item(46,11).Value = 6000
$ws.Cells.Item(46,12).Value = 0
$ws.Cells.Item(46,13).Value = -5844
$ws.Cells.Item(46,14).ClearContents()

$ws.Cells.Item(57,8).Value = 14815.25
$ws.Cells.Item(57,10).Value = 14815.25
$ws.Cells.Item(57,12).Value = 14815.25
$ws.Cells.Item(57,14).Value = -16455.25

$ws.Cells.Item(80,8).Value = 2553.8462
$ws.Cells.Item(80,9).Value = 2600
$ws.Cells.Item(80,10).Value = 2533.3333
$ws.Cells.Item(80,11).Value = 2600
$ws.Cells.Item(80,12).Value = 2533.3333
$ws.Cells.Item(80,13).Value = -1602
$ws.Cells.Item(80,14).Value = -4529.3333

$ws.Cells.Item(83,8).Value = 2553.8462
$ws.Cells.Item(83,9).Value = 2600
$ws.Cells.Item(83,10).Value = 2533.3333
$ws.Cells.Item(83,11).Value = 13000
$ws.Cells.Item(83,12).Value = 12666.6665
$ws.Cells.Item(83,13).Value = -8008
$ws.Cells.Item(83,14).Value = -22650.6665

$ws.Cells.Item(122,8).Value = 1458.0526
$ws.Cells.Item(122,9).Value = 807.1818
$ws.Cells.Item(122,10).Value = 2353
$ws.Cells.Item(122,11).Value = 2421.5454
$ws.Cells.Item(122,12).Value = 7059
$ws.Cells.Item(122,13).Value = 28.45460000000003
$ws.Cells.Item(122,14).Value = -11959

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132,8).Value = 10644861
$ws.Cells.Item(132,9).Value = 19240314
$ws.Cells.Item(132,10).Value = 2871.2856
$ws.Cells.Item(132,11).Value = 57720942
$ws.Cells.Item(132,12).Value = 8613.856800000001
$ws.Cells.Item(132,13).Value = -57718412
$ws.Cells.Item(132,14).Value = -13673.8568

$ws.Cells.Item(136,8).Value = 7150.304
$ws.Cells.Item(136,9).Value = 19286.715
$ws.Cells.Item(136,10).Value = 1840.625
$ws.Cells.Item(136,11).Value = 57860.145
$ws.Cells.Item(136,12).Value = 5521.875
$ws.Cells.Item(136,13).Value = -55310.145
$ws.Cells.Item(136,14).Value = -10621.875

$ws.Cells.Item(140,8).Value = 19350
$ws.Cells.Item(140,10).Value = 19350
$ws.Cells.Item(140,12).Value = 19350
$ws.Cells.Item(140,14).Value = -29710

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(3,8).Value = 0
$ws.Cells.Item(3,9).Value = 0
$ws.Cells.Item(3,11).Value = 0
$ws.Cells.Item(3,13).ClearContents()

$ws.Cells.Item(11,8).Value = 1016.6667
$ws.Cells.Item(11,10).Value = 1016.6667
$ws.Cells.Item(11,12).Value = 1016.6667
$ws.Cells.Item(11,14).Value = -1300.6667

$ws.Cells.Item(132,8).Value = 2133
$ws.Cells.Item(132,9).Value = 1427.7273
$ws.Cells.Item(132,10).Value = 2589.353
$ws.Cells.Item(132,11).Value = 4283.1819
$ws.Cells.Item(132,12).Value = 7768.059
$ws.Cells.Item(132,13).Value = -1753.1819
$ws.Cells.Item(132,14).Value = -12828.059

$ws.Cells.Item(136,8).Value = 1190.7805
$ws.Cells.Item(136,9).Value = 513.25
$ws.Cells.Item(136,10).Value = 2147.2942
$ws.Cells.Item(136,11).Value = 1539.75
$ws.Cells.Item(136,12).Value = 6441.882599999999
$ws.Cells.Item(136,13).Value = 1010.25
$ws.Cells.Item(136,14).Value = -11541.8826
